$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Harvest Area") values: 0.05 -> 0.01 for most rows, 0.05 -> 0 for rows 49-55
$rangesTo001 = @(
    @(2,6),
    @(9,18),
    @(22,30),
    @(34,42),
    @(45,48)
)
foreach ($r in $rangesTo001) {
    $ws.Range("C$($r[0]):C$($r[1])").Value = 0.01
}

$ws.Range("C49:C55").Value = 0

# Column E ("End Time") values: 20 -> 100 for all data rows (2-55)
$ws.Range("E2:E55").Value = 100

# Update the active selection shown in the sheet view to E5
$ws.Range("E5").Select()
